$wb = $excel.ActiveWorkbook

$wsTareas = $wb.Worksheets.Item("Tareas divididas")
$wsUsuarios = $wb.Worksheets.Item("Usuarios")

# Update B2 (Leandro), B3 (Ezequiel, new), B4 (Franco)
$wsTareas.Range("B2").Value = "Reportes, Armado procedimiento de suma del total a pagar, abm mesas, asignar mesero a mesa, Validaciones"
$wsTareas.Range("B3").Value = "Gestion de ordenes y pedidos, filtros y búsquedas mesero, Validaciones"
$wsTareas.Range("B4").Value = "Stock, Bajas logicas de categorias y menus, filtros y búsquedas gerencia, Validaciones"

# Column B width widened to fit new text
$wsTareas.Columns.Item(2).ColumnWidth = 99.85546875

# Selection / view changes
$wsTareas.Range("B2").Select()
$wsUsuarios.Range("B4").Select()

# Active sheet changes: Usuarios becomes the active/selected tab
$wsUsuarios.Activate()

Write-Host "done"
